$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-07-15 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-16 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("37-9=", $true, $false, $false, $false, $false, $true, 1, $false, "45+33=", 2) | Out-Null
$d.Content.Find.Execute("3+80=", $true, $false, $false, $false, $false, $true, 1, $false, "38+50=", 2) | Out-Null
$d.Content.Find.Execute("98+1=", $true, $false, $false, $false, $false, $true, 1, $false, "33-10=", 2) | Out-Null
$d.Content.Find.Execute("30+65=", $true, $false, $false, $false, $false, $true, 1, $false, "47-21=", 2) | Out-Null
$d.Content.Find.Execute("18+69=", $true, $false, $false, $false, $false, $true, 1, $false, "50+32=", 2) | Out-Null
$d.Content.Find.Execute("50+49=", $true, $false, $false, $false, $false, $true, 1, $false, "26+42=", 2) | Out-Null
$d.Content.Find.Execute("24-17=", $true, $false, $false, $false, $false, $true, 1, $false, "66+20=", 2) | Out-Null
$d.Content.Find.Execute("68+2=", $true, $false, $false, $false, $false, $true, 1, $false, "21+23=", 2) | Out-Null
$d.Content.Find.Execute("13+33=", $true, $false, $false, $false, $false, $true, 1, $false, "53-46=", 2) | Out-Null
$d.Content.Find.Execute("54+23=", $true, $false, $false, $false, $false, $true, 1, $false, "16-8=", 2) | Out-Null
$d.Content.Find.Execute("20+52=", $true, $false, $false, $false, $false, $true, 1, $false, "20-17=", 2) | Out-Null
$d.Content.Find.Execute("36+16=", $true, $false, $false, $false, $false, $true, 1, $false, "84-15=", 2) | Out-Null
$d.Content.Find.Execute("93-75=", $true, $false, $false, $false, $false, $true, 1, $false, "46-42=", 2) | Out-Null
$d.Content.Find.Execute("98-1=", $true, $false, $false, $false, $false, $true, 1, $false, "54-22=", 2) | Out-Null
$d.Content.Find.Execute("38-15=", $true, $false, $false, $false, $false, $true, 1, $false, "55+18=", 2) | Out-Null
$d.Content.Find.Execute("97-85=", $true, $false, $false, $false, $false, $true, 1, $false, "69+9=", 2) | Out-Null
$d.Content.Find.Execute("2+60=", $true, $false, $false, $false, $false, $true, 1, $false, "7+78=", 2) | Out-Null
$d.Content.Find.Execute("55+20=", $true, $false, $false, $false, $false, $true, 1, $false, "55+21=", 2) | Out-Null
$d.Content.Find.Execute("86-62=", $true, $false, $false, $false, $false, $true, 1, $false, "57-18=", 2) | Out-Null
$d.Content.Find.Execute("73+24=", $true, $false, $false, $false, $false, $true, 1, $false, "11-11=", 2) | Out-Null
$d.Content.Find.Execute("12+31=", $true, $false, $false, $false, $false, $true, 1, $false, "14+61=", 2) | Out-Null
$d.Content.Find.Execute("71+7=", $true, $false, $false, $false, $false, $true, 1, $false, "20+74=", 2) | Out-Null
$d.Content.Find.Execute("77-46=", $true, $false, $false, $false, $false, $true, 1, $false, "12+15=", 2) | Out-Null
$d.Content.Find.Execute("17+62=", $true, $false, $false, $false, $false, $true, 1, $false, "57+31=", 2) | Out-Null
$d.Content.Find.Execute("37-29=", $true, $false, $false, $false, $false, $true, 1, $false, "83-0=", 2) | Out-Null
$d.Content.Find.Execute("58-18=", $true, $false, $false, $false, $false, $true, 1, $false, "65-29=", 2) | Out-Null
$d.Content.Find.Execute("95-27=", $true, $false, $false, $false, $false, $true, 1, $false, "31+17=", 2) | Out-Null
$d.Content.Find.Execute("2+59=", $true, $false, $false, $false, $false, $true, 1, $false, "22+1=", 2) | Out-Null
$d.Content.Find.Execute("27+58=", $true, $false, $false, $false, $false, $true, 1, $false, "47+6=", 2) | Out-Null
$d.Content.Find.Execute("24+31=", $true, $false, $false, $false, $false, $true, 1, $false, "52-40=", 2) | Out-Null
$d.Content.Find.Execute("75-48=", $true, $false, $false, $false, $false, $true, 1, $false, "13+6=", 2) | Out-Null
$d.Content.Find.Execute("93+5=", $true, $false, $false, $false, $false, $true, 1, $false, "56+7=", 2) | Out-Null
$d.Content.Find.Execute("65-47=", $true, $false, $false, $false, $false, $true, 1, $false, "80-35=", 2) | Out-Null
$d.Content.Find.Execute("89-58=", $true, $false, $false, $false, $false, $true, 1, $false, "18+15=", 2) | Out-Null
$d.Content.Find.Execute("74-68=", $true, $false, $false, $false, $false, $true, 1, $false, "60+1=", 2) | Out-Null
$d.Content.Find.Execute("28+28=", $true, $false, $false, $false, $false, $true, 1, $false, "34+49=", 2) | Out-Null
$d.Content.Find.Execute("40+39=", $true, $false, $false, $false, $false, $true, 1, $false, "42+41=", 2) | Out-Null
$d.Content.Find.Execute("47-29=", $true, $false, $false, $false, $false, $true, 1, $false, "93-4=", 2) | Out-Null
$d.Content.Find.Execute("8-5=", $true, $false, $false, $false, $false, $true, 1, $false, "20+58=", 2) | Out-Null
$d.Content.Find.Execute("58-39=", $true, $false, $false, $false, $false, $true, 1, $false, "72+9=", 2) | Out-Null
$d.Content.Find.Execute("57+12=", $true, $false, $false, $false, $false, $true, 1, $false, "3+47=", 2) | Out-Null
$d.Content.Find.Execute("41+10=", $true, $false, $false, $false, $false, $true, 1, $false, "89-5=", 2) | Out-Null
$d.Content.Find.Execute("84+2=", $true, $false, $false, $false, $false, $true, 1, $false, "31+25=", 2) | Out-Null
$d.Content.Find.Execute("29+16=", $true, $false, $false, $false, $false, $true, 1, $false, "47-17=", 2) | Out-Null
$d.Content.Find.Execute("22+7=", $true, $false, $false, $false, $false, $true, 1, $false, "78-23=", 2) | Out-Null
$d.Content.Find.Execute("35+63=", $true, $false, $false, $false, $false, $true, 1, $false, "64+29=", 2) | Out-Null
$d.Content.Find.Execute("7+11=", $true, $false, $false, $false, $false, $true, 1, $false, "63+13=", 2) | Out-Null
$d.Content.Find.Execute("43+29=", $true, $false, $false, $false, $false, $true, 1, $false, "16+35=", 2) | Out-Null
$d.Content.Find.Execute("98-21=", $true, $false, $false, $false, $false, $true, 1, $false, "60-53=", 2) | Out-Null
$d.Content.Find.Execute("90-5=", $true, $false, $false, $false, $false, $true, 1, $false, "79-75=", 2) | Out-Null
$d.Content.Find.Execute("27-2=", $true, $false, $false, $false, $false, $true, 1, $false, "53+43=", 2) | Out-Null
$d.Content.Find.Execute("36-8=", $true, $false, $false, $false, $false, $true, 1, $false, "45+37=", 2) | Out-Null
$d.Content.Find.Execute("42+25=", $true, $false, $false, $false, $false, $true, 1, $false, "47-37=", 2) | Out-Null
$d.Content.Find.Execute("98-3=", $true, $false, $false, $false, $false, $true, 1, $false, "91-78=", 2) | Out-Null
$d.Content.Find.Execute("36+27=", $true, $false, $false, $false, $false, $true, 1, $false, "61+15=", 2) | Out-Null
$d.Content.Find.Execute("19+53=", $true, $false, $false, $false, $false, $true, 1, $false, "12+23=", 2) | Out-Null
$d.Content.Find.Execute("10+78=", $true, $false, $false, $false, $false, $true, 1, $false, "40+19=", 2) | Out-Null
$d.Content.Find.Execute("44+11=", $true, $false, $false, $false, $false, $true, 1, $false, "79-53=", 2) | Out-Null
$d.Content.Find.Execute("83-59=", $true, $false, $false, $false, $false, $true, 1, $false, "39-19=", 2) | Out-Null
$d.Content.Find.Execute("64-41=", $true, $false, $false, $false, $false, $true, 1, $false, "13-0=", 2) | Out-Null
$d.Content.Find.Execute("56-1=", $true, $false, $false, $false, $false, $true, 1, $false, "96-44=", 2) | Out-Null
$d.Content.Find.Execute("51-29=", $true, $false, $false, $false, $false, $true, 1, $false, "55+24=", 2) | Out-Null
$d.Content.Find.Execute("55-46=", $true, $false, $false, $false, $false, $true, 1, $false, "79-74=", 2) | Out-Null
$d.Content.Find.Execute("96-45=", $true, $false, $false, $false, $false, $true, 1, $false, "8+50=", 2) | Out-Null
$d.Content.Find.Execute("58-20=", $true, $false, $false, $false, $false, $true, 1, $false, "15+77=", 2) | Out-Null
$d.Content.Find.Execute("17+25=", $true, $false, $false, $false, $false, $true, 1, $false, "72-22=", 2) | Out-Null
$d.Content.Find.Execute("93+0=", $true, $false, $false, $false, $false, $true, 1, $false, "21-8=", 2) | Out-Null
$d.Content.Find.Execute("25+29=", $true, $false, $false, $false, $false, $true, 1, $false, "59-19=", 2) | Out-Null
$d.Content.Find.Execute("68-43=", $true, $false, $false, $false, $false, $true, 1, $false, "99-73=", 2) | Out-Null
$d.Content.Find.Execute("19+79=", $true, $false, $false, $false, $false, $true, 1, $false, "29+21=", 2) | Out-Null
$d.Content.Find.Execute("27+27=", $true, $false, $false, $false, $false, $true, 1, $false, "38+45=", 2) | Out-Null
$d.Content.Find.Execute("66+29=", $true, $false, $false, $false, $false, $true, 1, $false, "16+13=", 2) | Out-Null
$d.Content.Find.Execute("39-37=", $true, $false, $false, $false, $false, $true, 1, $false, "50-37=", 2) | Out-Null
$d.Content.Find.Execute("55+2=", $true, $false, $false, $false, $false, $true, 1, $false, "64-56=", 2) | Out-Null
$d.Content.Find.Execute("58-42=", $true, $false, $false, $false, $false, $true, 1, $false, "38+17=", 2) | Out-Null
$d.Content.Find.Execute("71-36=", $true, $false, $false, $false, $false, $true, 1, $false, "52-45=", 2) | Out-Null
$d.Content.Find.Execute("69-32=", $true, $false, $false, $false, $false, $true, 1, $false, "57+41=", 2) | Out-Null
$d.Content.Find.Execute("52-46=", $true, $false, $false, $false, $false, $true, 1, $false, "5+82=", 2) | Out-Null
$d.Content.Find.Execute("78-15=", $true, $false, $false, $false, $false, $true, 1, $false, "96-46=", 2) | Out-Null
$d.Content.Find.Execute("93-77=", $true, $false, $false, $false, $false, $true, 1, $false, "83-52=", 2) | Out-Null
$d.Content.Find.Execute("60+24=", $true, $false, $false, $false, $false, $true, 1, $false, "71-37=", 2) | Out-Null
$d.Content.Find.Execute("11+87=", $true, $false, $false, $false, $false, $true, 1, $false, "2+97=", 2) | Out-Null
$d.Content.Find.Execute("69+3=", $true, $false, $false, $false, $false, $true, 1, $false, "58+14=", 2) | Out-Null
$d.Content.Find.Execute("67-56=", $true, $false, $false, $false, $false, $true, 1, $false, "93-43=", 2) | Out-Null
$d.Content.Find.Execute("96-17=", $true, $false, $false, $false, $false, $true, 1, $false, "79-49=", 2) | Out-Null
$d.Content.Find.Execute("21-4=", $true, $false, $false, $false, $false, $true, 1, $false, "62-34=", 2) | Out-Null
$d.Content.Find.Execute("19+38=", $true, $false, $false, $false, $false, $true, 1, $false, "92-63=", 2) | Out-Null
$d.Content.Find.Execute("58-44=", $true, $false, $false, $false, $false, $true, 1, $false, "72-9=", 2) | Out-Null
$d.Content.Find.Execute("56-41=", $true, $false, $false, $false, $false, $true, 1, $false, "92-1=", 2) | Out-Null
$d.Content.Find.Execute("16-7=", $true, $false, $false, $false, $false, $true, 1, $false, "43-18=", 2) | Out-Null
$d.Content.Find.Execute("67+6=", $true, $false, $false, $false, $false, $true, 1, $false, "42+5=", 2) | Out-Null
$d.Content.Find.Execute("24+12=", $true, $false, $false, $false, $false, $true, 1, $false, "26+18=", 2) | Out-Null
$d.Content.Find.Execute("17+41=", $true, $false, $false, $false, $false, $true, 1, $false, "74-16=", 2) | Out-Null
$d.Content.Find.Execute("80+5=", $true, $false, $false, $false, $false, $true, 1, $false, "69-45=", 2) | Out-Null
$d.Content.Find.Execute("65+32=", $true, $false, $false, $false, $false, $true, 1, $false, "12+81=", 2) | Out-Null
$d.Content.Find.Execute("44+43=", $true, $false, $false, $false, $false, $true, 1, $false, "12+30=", 2) | Out-Null
$d.Content.Find.Execute("9+69=", $true, $false, $false, $false, $false, $true, 1, $false, "16+50=", 2) | Out-Null
$d.Content.Find.Execute("73+21=", $true, $false, $false, $false, $false, $true, 1, $false, "71-23=", 2) | Out-Null
$d.Content.Find.Execute("82+7=", $true, $false, $false, $false, $false, $true, 1, $false, "35+28=", 2) | Out-Null
$d.Content.Find.Execute("29+66=", $true, $false, $false, $false, $false, $true, 1, $false, "80-23=", 2) | Out-Null
